# Updated cryptos list on Fri Jul  7 08:20:15 UTC 2023 with GitHub Actions
# Refreshes the price/volume columns (D, E) for each coin row, and swaps the
# RenderToken/Quant rows (43 & 45) back into their refreshed ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NumberFormat "@" forces text parsing so numeric-looking strings like
# "1.000" or "233.45" are stored verbatim instead of being coerced into
# Excel numbers; resetting the Style back to "Normal" afterwards strips the
# temporary text format so no stray cell style is left behind.

$ws.Range("D2").Value = "30.070.82"
$ws.Range("E2").Value = "  -2.46%  "

$ws.Range("D3").Value = "1.860.21"
$ws.Range("E3").Value = "  -3.55%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4652"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.67%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2799"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.81%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06536"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.75%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07804"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.47"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.22%  "

$ws.Range("D13").Value = "1.872.57"
$ws.Range("E13").Value = "  -2.94%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.097"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6633"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "281.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.55%  "

$ws.Range("D17").Value = "30.108.51"
$ws.Range("E17").Value = "  -2.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.459"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.51%  "

$ws.Range("D21").Value = "2.105.98"
$ws.Range("E21").Value = "  -3.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007205"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.113"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.307"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.68%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.904"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.333"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09540"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.419"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.466"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.078"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04627"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.096"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6979"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.694"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01845"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.266"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.503"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8531"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.56%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.79%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.0000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.902"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4138"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.80%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "988.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.35%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.172"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.249"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1133"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.13%  "
